# Auto-generated script to apply odds updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 17
$ws.Range("I17").Value = 2.42
$ws.Range("K17").Value = 2.12
$ws.Range("L17").Value = 3
$ws.Range("N17").Value = 7.6
$ws.Range("O17").Value = 1.27
$ws.Range("P17").Value = 3.4
$ws.Range("Q17").Value = 1.83
$ws.Range("R17").Value = 1.91
$ws.Range("S17").Value = 1.39
$ws.Range("T17").Value = 2.77
$ws.Range("U17").Value = 1.65
$ws.Range("V17").Value = 2.1
$ws.Range("W17").Value = 9.5
$ws.Range("X17").Value = 14
$ws.Range("AA17").Value = 22
$ws.Range("AB17").Value = 29
$ws.Range("AC17").Value = 7.6
$ws.Range("AD17").Value = 6.5
$ws.Range("AE17").Value = 13
$ws.Range("AF17").Value = 55
$ws.Range("AG17").Value = 400
$ws.Range("AH17").Value = 9
$ws.Range("AL17").Value = 19
$ws.Range("AM17").Value = 27
$ws.Range("AN17").Value = 4.65
$ws.Range("AP17").Value = 22
$ws.Range("AQ17").Value = 65
$ws.Range("AR17").Value = 100
$ws.Range("AT17").Value = 2.77
$ws.Range("AU17").Value = 7
$ws.Range("AV17").Value = 60
$ws.Range("AW17").Value = 4.4
$ws.Range("AX17").Value = 13
$ws.Range("AY17").Value = 20
$ws.Range("AZ17").Value = 55

# Row 21
$ws.Range("G21").Value = 2.1
$ws.Range("H21").Value = 3.25
$ws.Range("I21").Value = 3.7
$ws.Range("L21").Value = 4.5
$ws.Range("W21").Value = 6
$ws.Range("Z21").Value = 19
$ws.Range("AF21").Value = 67
$ws.Range("AH21").Value = 8
$ws.Range("AL21").Value = 34
$ws.Range("AN21").Value = 4
$ws.Range("AR21").Value = 81

# Row 22
$ws.Range("G22").Value = 1.4
$ws.Range("I22").Value = 7.5
$ws.Range("N22").Value = 12
$ws.Range("U22").Value = 2.05
$ws.Range("V22").Value = 1.7
$ws.Range("AC22").Value = 12
$ws.Range("AG22").Value = 451
$ws.Range("AI22").Value = 41
$ws.Range("AN22").Value = 3.25
$ws.Range("AO22").Value = 6.5
$ws.Range("AQ22").Value = 19
$ws.Range("AS22").Value = 151
$ws.Range("AU22").Value = 9.5
$ws.Range("AV22").Value = 67
$ws.Range("BB22").Value = 351

# Row 24
$ws.Range("G24").Value = 2.25
$ws.Range("I24").Value = 3.2
$ws.Range("J24").Value = 3.2
$ws.Range("M24").Value = 1.13
$ws.Range("N24").Value = 6
$ws.Range("Y24").Value = 10
$ws.Range("Z24").Value = 21
$ws.Range("AK24").Value = 41
$ws.Range("AN24").Value = 4
$ws.Range("AR24").Value = 81
$ws.Range("AY24").Value = 41
$ws.Range("AZ24").Value = 81

# Row 28
$ws.Range("G28").Value = 2.05
$ws.Range("I28").Value = 3.2
$ws.Range("J28").Value = 2.63
$ws.Range("L28").Value = 3.6
$ws.Range("N28").Value = 17
$ws.Range("Q28").Value = 1.6
$ws.Range("R28").Value = 2.3
$ws.Range("S28").Value = 1.29
$ws.Range("T28").Value = 3.5
$ws.Range("W28").Value = 11
$ws.Range("X28").Value = 12
$ws.Range("AC28").Value = 17
$ws.Range("AF28").Value = 34
$ws.Range("AK28").Value = 34
$ws.Range("AO28").Value = 11
$ws.Range("AT28").Value = 3.5
$ws.Range("BA28").Value = 51

# Row 30
$ws.Range("N30").Value = 19

# Row 41
$ws.Range("O41").Value = 1.3
$ws.Range("P41").Value = 3.5
$ws.Range("Q41").Value = 2
$ws.Range("R41").Value = 1.85

# Row 46
$ws.Range("O46").Value = 1.33
$ws.Range("P46").Value = 3.25
$ws.Range("R46").Value = 1.75
$ws.Range("S46").Value = 1.44
$ws.Range("T46").Value = 2.63
$ws.Range("AT46").Value = 2.63

# Row 47
$ws.Range("O47").Value = 1.3
$ws.Range("P47").Value = 3.4
$ws.Range("Q47").Value = 2.03
$ws.Range("R47").Value = 1.83

# Row 50
$ws.Range("I50").Value = 2.63
$ws.Range("N50").Value = 12
$ws.Range("Y50").Value = 10
$ws.Range("AD50").Value = 7
$ws.Range("AK50").Value = 26
$ws.Range("AN50").Value = 4.75
$ws.Range("AZ50").Value = 41

# Row 52
$ws.Range("G52").Value = 5.75
$ws.Range("H52").Value = 3.9
$ws.Range("I52").Value = 1.5
$ws.Range("K52").Value = 2.25
$ws.Range("L52").Value = 2.1
$ws.Range("Y52").Value = 19
$ws.Range("AA52").Value = 51
$ws.Range("AD52").Value = 8
$ws.Range("AN52").Value = 7.5
$ws.Range("AO52").Value = 34
$ws.Range("AU52").Value = 9
$ws.Range("AV52").Value = 67
$ws.Range("AW52").Value = 3.4
$ws.Range("AZ52").Value = 23

# Row 54
$ws.Range("G54").Value = 2.2
$ws.Range("I54").Value = 3.1
$ws.Range("L54").Value = 3.75
$ws.Range("N54").Value = 9
$ws.Range("AZ54").Value = 51

# Row 55
$ws.Range("G55").Value = 2.15
$ws.Range("J55").Value = 2.88
$ws.Range("L55").Value = 4
$ws.Range("AS55").Value = 201

# Row 56
$ws.Range("M56").Value = 1.03
$ws.Range("O56").Value = 1.17

# Row 57
$ws.Range("M57").Value = 1.03
$ws.Range("O57").Value = 1.19

# Row 60
$ws.Range("M60").Value = 1.01
$ws.Range("O60").Value = 1.12

# Row 76
$ws.Range("L76").Value = 3.25
$ws.Range("AO76").Value = 19
$ws.Range("AS76").Value = 301

# Row 77
$ws.Range("M77").Value = 1.03
$ws.Range("O77").Value = 1.19

# Row 78
$ws.Range("M78").Value = 1.05
$ws.Range("O78").Value = 1.33

# Row 81
$ws.Range("M81").Value = 1.05
$ws.Range("O81").Value = 1.37

# Row 82
$ws.Range("K82").Value = 1.95
$ws.Range("M82").Value = 1.07
$ws.Range("O82").Value = 1.47

# Row 83
$ws.Range("M83").Value = 1.03
$ws.Range("O83").Value = 1.19

# Row 90
$ws.Range("I90").Value = 3.5
$ws.Range("L90").Value = 3.75
$ws.Range("N90").Value = 12
$ws.Range("O90").Value = 1.22
$ws.Range("P90").Value = 4
$ws.Range("Q90").Value = 1.8
$ws.Range("R90").Value = 2
$ws.Range("X90").Value = 11
$ws.Range("Y90").Value = 9
$ws.Range("Z90").Value = 19
$ws.Range("AD90").Value = 7
$ws.Range("AL90").Value = 26
$ws.Range("AM90").Value = 29
$ws.Range("AN90").Value = 4.33
$ws.Range("AY90").Value = 23
$ws.Range("BA90").Value = 67

# Row 91
$ws.Range("G91").Value = 1.7
$ws.Range("H91").Value = 3.4
$ws.Range("I91").Value = 5.75
$ws.Range("M91").Value = 1.08
$ws.Range("N91").Value = 8
$ws.Range("X91").Value = 7
$ws.Range("Z91").Value = 12
$ws.Range("AJ91").Value = 19
$ws.Range("AN91").Value = 3.5
